# "Almost forgot to add Ancient Gear Token"
#
# SR03-JP (Structure Deck R: Machine Dragon Re-Volt) is missing the token
# that belongs with the structure deck — add it as the next row, then
# leave the workbook's view state the way the author left it: cursor
# parked past the data on SR03-JP (H5) and the SPDS-JP sheet scrolled/
# selected down at A23 (and still the active tab).

$wb = $excel.ActiveWorkbook

# --- SR03-JP: add the missing "Ancient Gear Token" row ------------------
$ws4 = $wb.Worksheets.Item("SR03-JP")

$ws4.Range("A5").Value = "Ancient Gear Token"
$ws4.Range("B5").Value = 100303121
$ws4.Range("C5").Value = ":"
$ws4.Range("E5").Value = ";"

# Leave the cursor where the author's last keystroke landed (one column
# past the used range) without disturbing which tab is active overall.
[void]$ws4.Range("H5").Select()

# --- SPDS-JP: restore it as the active sheet, scrolled/selected at A23 --
$ws3 = $wb.Worksheets.Item("SPDS-JP")
$ws3.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1

[void]$ws3.Range("A23").Select()
